$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape "Rectangle 4" (the first "Trial parameters file" / PARAMETER_TRIAL box)
# becomes the "Local parameters file" / LOCALPARAM_INFO box.
# ---------------------------------------------------------------------------
$shp1 = $s.Shapes.Item(1)
$tr1 = $shp1.TextFrame.TextRange
$full1 = $tr1.Text

$old1a = "Trial parameters file"
$idx1a = $full1.IndexOf($old1a)
$c1a = $tr1.Characters($idx1a + 1, ("Trial parameters ").Length)
$c1a.Text = "Local parameters "

# refresh text after the first edit before locating the second target
$full1b = $tr1.Text
$old1b = "PARAMETER_TRIAL "
$idx1b = $full1b.IndexOf($old1b)
$c1b = $tr1.Characters($idx1b + 1, $old1b.Length)
$c1b.Text = "LOCALPARAM_INFO "

# ---------------------------------------------------------------------------
# Shape "Rectangle 9" (the second "Trial parameters file " / PARAMETER_TRIAL
# box) gets its runs re-split without changing the visible text.
# ---------------------------------------------------------------------------
$shp2 = $s.Shapes.Item(5)
$tr2 = $shp2.TextFrame.TextRange
$full2 = $tr2.Text

$old2a = "Trial parameters "
$idx2a = $full2.IndexOf($old2a)
$c2a = $tr2.Characters($idx2a + 1, $old2a.Length)
$c2a.Text = "Trial parameters "

$full2b = $tr2.Text
$old2b = "PARAMETER_TRIAL -- "
$idx2b = $full2b.IndexOf($old2b)
$c2b = $tr2.Characters($idx2b + 1, $old2b.Length)
$c2b.Text = "PARAMETER_TRIAL -- "

$full2c = $tr2.Text
$old2c = "ASCII"
$idx2c = $full2c.LastIndexOf($old2c)
$c2c = $tr2.Characters($idx2c + 1, $old2c.Length)
$c2c.Text = "NC"

# ---------------------------------------------------------------------------
# Ovals "3" and "4": drop the trailing endParaRPr element. Deleting the
# existing run and retyping the same digit reproduces the text without
# PowerPoint re-adding the (now redundant) end paragraph run properties.
# ---------------------------------------------------------------------------
$shp3 = $s.Shapes.Item(10)
$tr3 = $shp3.TextFrame.TextRange
$text3 = $tr3.Text
[void]$tr3.Delete()
$tr3.Text = $text3

$shp4 = $s.Shapes.Item(11)
$tr4 = $shp4.TextFrame.TextRange
$text4 = $tr4.Text
[void]$tr4.Delete()
$tr4.Text = $text4
